# Update stock prices (股价) for 2023-06, sheets 个人持仓 and 家庭持仓
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C2").Value = 46.4
$ws1.Range("C3").Value = 44.89
$ws1.Range("C4").Value = 53.45
$ws1.Range("C5").Value = 52.48
$ws1.Range("C6").Value = 31.88
$ws1.Range("C7").Value = 43.37
$ws1.Range("C8").Value = 27.22
$ws1.Range("C9").Value = 27.04
$ws1.Range("C10").Value = 27.15
$ws1.Range("C11").Value = 131.35
$ws1.Range("C12").Value = 163.57
$ws1.Range("C13").Value = 209.57
$ws1.Range("C14").Value = 0.764
$ws1.Range("C15").Value = 10.86
$ws1.Range("C16").Value = 28.32
$ws1.Range("C17").Value = 24.49
$ws1.Range("C18").Value = 18.31
$ws1.Range("C19").Value = 36.21
$ws1.Range("C20").Value = 36.51
$ws1.Range("C21").Value = 25.29
$ws1.Range("C22").Value = 110.53
$ws1.Range("C23").Value = 3.882
$ws1.Range("C24").Value = 4.349

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C2").Value = 46.4
$ws2.Range("C3").Value = 44.89
$ws2.Range("C4").Value = 53.45
$ws2.Range("C5").Value = 52.48
$ws2.Range("C6").Value = 31.88
$ws2.Range("C7").Value = 43.37
$ws2.Range("C8").Value = 27.22
$ws2.Range("C9").Value = 27.04
$ws2.Range("C10").Value = 27.15
$ws2.Range("C11").Value = 131.35
$ws2.Range("C12").Value = 209.57
$ws2.Range("C13").Value = 163.57
$ws2.Range("C14").Value = 185.07
$ws2.Range("C15").Value = 0.764
$ws2.Range("C16").Value = 46.85
$ws2.Range("C17").Value = 10.86
$ws2.Range("C18").Value = 28.32
$ws2.Range("C19").Value = 24.49
$ws2.Range("C20").Value = 18.31
$ws2.Range("C21").Value = 36.21
$ws2.Range("C22").Value = 36.51
$ws2.Range("C23").Value = 25.29
$ws2.Range("C24").Value = 110.53
$ws2.Range("C25").Value = 3.882
$ws2.Range("C26").Value = 1.118
$ws2.Range("C27").Value = 4.349

$wb.Save()
